$d = $word.ActiveDocument

# 1. Split the "Critério" run: insert a manual line break after "em que:"
$search1 = "Para compor a Nota Final (NF) serão consideradas as avaliações de Provas Escritas (P1 e P2) e Exercícios Propostos e Seminários (T) em que:NF = (P1 + P2 + T)/3.  Serão considerados aprovados os alunos que obtiverem: NF maior ou igual a 5,0. Serão considerados reprovados os alunos que obtiverem: NS menor que 3,0 Para os alunos em que NS é maior ou igual a 3,0 e menor que 5,0 será dada uma prova de recuperação (R)."
$replace1 = "Para compor a Nota Final (NF) serão consideradas as avaliações de Provas Escritas (P1 e P2) e Exercícios Propostos e Seminários (T) em que:^lNF = (P1 + P2 + T)/3.  Serão considerados aprovados os alunos que obtiverem: NF maior ou igual a 5,0. Serão considerados reprovados os alunos que obtiverem: NS menor que 3,0 Para os alunos em que NS é maior ou igual a 3,0 e menor que 5,0 será dada uma prova de recuperação (R)."
$d.Content.Find.Execute($search1, $true, $false, $false, $false, $false, $true, 1, $false, $replace1, 2)

# 2. Split the Bibliografia run into five references separated by manual line breaks
$search2 = "1. F.P. BEER, E.R. JOHNSTON, J.T. DeWOLF, D.MAZUREK. Estática e Mecânica dos Materiais. São Paulo: McGraw Hill, 2013, 728p.2. F.P. BEER, E.R. JOHNSTON, E. RUSSEL. Mecânica vetorial para engenheiros: Estática. São Paulo: McGraw Hill. 9a Ed., 2012, 626p. 3. HIBBELER, R.C. Mecânica para engenharia vol.1: estática. São Paulo: Pearson Prentice Hall, 12a Ed., 2011.4. MERIAM, J.L. KRAIGE, L.G. Mecânica para engenharia – Estática. Grupo GEN Editora LTC, 6a Ed., 2009, 364p. 5. RUIZ, C.C.de La P. Fundamentos de mecânica para engenharia – Estática. Grupo GEN Editora LTC, 2017, 306p."
$replace2 = "1. F.P. BEER, E.R. JOHNSTON, J.T. DeWOLF, D.MAZUREK. Estática e Mecânica dos Materiais. São Paulo: McGraw Hill, 2013, 728p.^l2. F.P. BEER, E.R. JOHNSTON, E. RUSSEL. Mecânica vetorial para engenheiros: Estática. São Paulo: McGraw Hill. 9a Ed., 2012, 626p. ^l3. HIBBELER, R.C. Mecânica para engenharia vol.1: estática. São Paulo: Pearson Prentice Hall, 12a Ed., 2011.^l4. MERIAM, J.L. KRAIGE, L.G. Mecânica para engenharia – Estática. Grupo GEN Editora LTC, 6a Ed., 2009, 364p. ^l5. RUIZ, C.C.de La P. Fundamentos de mecânica para engenharia – Estática. Grupo GEN Editora LTC, 2017, 306p."
$d.Content.Find.Execute($search2, $true, $false, $false, $false, $false, $true, 1, $false, $replace2, 2)
